$d = $word.ActiveDocument

# The document ends with two empty paragraphs right before the section
# break. The second (last) of those gets the new "start of term paper"
# text, with a first-line indent of 0.25" (360 twips / 18 pt).
$p = $d.Paragraphs.Last
$p.Range.Text = "The data set includes information about the attrition rate for employees within the healthcare field. The meaning of employee attrition is the departure of employees from the organization for any reason whether that be voluntary or involuntary, including resignation, termination, death, or retirement. Companies to avoid attrition rates being too high is to replace those who are either leaving voluntarily or involuntary. The data set should provide insights into whether a company in the healthcare field was replacing their employees that were leaving the field, or if they continued to have a gradual but deliberate reduction in staff for any reason."
$p.Range.ParagraphFormat.FirstLineIndent = 18
